$wb = $excel.ActiveWorkbook

# Rename the existing sheet and update its data/selection
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "addCustomerFlow"

$ws1.Range("A1").Value = "FirstName"
$ws1.Range("B1").Value = "LastName"
$ws1.Range("C1").Value = "PostCode"
$ws1.Range("D1").Value = "SuccessMessage"

$ws1.Range("A2").Value = "Jack"
$ws1.Range("B2").Value = "Daniel"
$ws1.Range("C2").Value = "JD12345"
$ws1.Range("D2").Value = "Customer added successfully"

$ws1.Range("B3").Select()

# Add the new sheet after the first one
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "openAccountFlow"

$ws2.Range("A1").Value = "Customer"
$ws2.Range("B1").Value = "Currency"
$ws2.Range("C1").Value = "SuccessMessage"

$ws2.Range("A2").Value = "Harry Potter"
$ws2.Range("B2").Value = "Rupee"
$ws2.Range("C2").Value = "Account created successfully"

$ws2.Range("B6").Select()

$null = $ws2.Columns.Item(1).AutoFit()
$null = $ws2.Columns.Item(3).AutoFit()
